# Restore revision change: update cell C10 on the active sheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
